# memberType and user/companies import done
#
# - Removes the obsolete "Admin" column (K) entirely, shifting every
#   column to its right one position to the left (L->K, M->L, ... Q->P).
# - "joia paga" (column J) switches from a free-text marker ("v"/"f") to
#   a real boolean: rows 2 and 3 become TRUE.
# - Row 2's "date de entrada no clube" (I2) is updated to 2020-05-09
#   (serial 43960).
# - Selection moves to I14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "date de entrada no clube" for the first user row.
$ws.Range("I2").Value = 43960

# "joia paga" becomes a real boolean column; rows 2 & 3 are paid (TRUE).
$ws.Range("J2").Value = $true
$ws.Range("J3").Value = $true

# Drop the whole "Admin" column (K) - everything to the right shifts left.
$ws.Range("K1").EntireColumn.Delete()

# Match the saved selection from the edit.
[void]$ws.Range("I14").Select()
